$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shorten transmitter key value
$ws.Range("B2").Value = "a35df02ffc4b9fadd1b3"

# Rows 4-6: actions table - replace hex-only codes with "letter (hex)" codes
$ws.Range("B4").Value = "p (0x70)"
$ws.Range("B5").Value = "r (0x72)"
$ws.Range("B6").Value = "b (0x62)"

# Rows 8-11: avr-->server table - shift values up (drop old row 8 "ready for
# transmission" entry) and replace hex-only codes with "letter (hex)" codes
$ws.Range("B8").Value = "b (0x62)"
$ws.Range("C8").Value = "battery level transmission (2 bytes)"

$ws.Range("B9").Value = "a (0x61)"
$ws.Range("C9").Value = "bad action"

$ws.Range("B10").Value = "k (0x6b)"
$ws.Range("C10").Value = "bad key"

$ws.Range("B11").Value = "o (0x6f)"
$ws.Range("C11").Value = "OK"

# Row 12 no longer exists - clear its former contents
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

# Update the active selection
$ws.Range("B8").Select()
